$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status refresh: "Ready for handoff" -> "In Translation" everywhere it appears ---
# Overview sheet: zh-cn / de-de status columns (E & F), rows 2-3
$wsOverview.Range("E2:F3").Value = "In Translation"

# Per-locale sheets: Status column (C), rows 2-3
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- Narrow the status-date columns that used to be sized for "Ready for handoff" ---
# Overview: columns E and F (zh-cn / de-de)
$wsOverview.Range("E:F").ColumnWidth = 12.42

# Per-locale sheets: Status column (C)
$wsZhCn.Range("C:C").ColumnWidth = 12.42
$wsDeDe.Range("C:C").ColumnWidth = 12.42
